$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 24964.363
$ws.Range("I51").Value = 29187.375
$ws.Range("K51").Value = 29187.375
$ws.Range("M51").Value = -28703.375
$ws.Range("H69").Value = 76753.766
$ws.Range("I69").Value = 11466.5
$ws.Range("K69").Value = 34399.5
$ws.Range("M69").Value = -33525.5
$ws.Range("H72").Value = 76753.766
$ws.Range("I72").Value = 11466.5
$ws.Range("K72").Value = 103198.5
$ws.Range("M72").Value = -98830.5
$ws.Range("H80").Value = 471.2
$ws.Range("I80").Value = 449.2857
$ws.Range("J80").Value = 522.3333
$ws.Range("K80").Value = 1347.8571
$ws.Range("L80").Value = 1566.9999
$ws.Range("M80").Value = -349.8571000000002
$ws.Range("N80").Value = -3562.9999
$ws.Range("H83").Value = 471.2
$ws.Range("I83").Value = 449.2857
$ws.Range("J83").Value = 522.3333
$ws.Range("K83").Value = 4043.5713
$ws.Range("L83").Value = 4700.9997
$ws.Range("M83").Value = 948.4286999999999
$ws.Range("N83").Value = -14684.9997
$ws.Range("H88").Value = 2333
$ws.Range("I88").Value = 2499.5
$ws.Range("K88").Value = 2499.5
$ws.Range("M88").Value = -2093.5
$ws.Range("H91").Value = 2333
$ws.Range("I91").Value = 2499.5
$ws.Range("K91").Value = 2499.5
$ws.Range("M91").Value = -1095.5
$ws.Range("H116").Value = 13558.333
$ws.Range("I116").Value = 21224
$ws.Range("K116").Value = 21224
$ws.Range("M116").Value = -17782
$ws.Range("H132").Value = 3610.7307
$ws.Range("I132").Value = 3515.16
$ws.Range("K132").Value = 10545.48
$ws.Range("M132").Value = -8015.48
$ws.Range("H141").Value = 5006.75
$ws.Range("I141").Value = 4054.0527
$ws.Range("K141").Value = 12162.1581
$ws.Range("M141").Value = -6982.158100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 64998.75
$ws.Range("J44").Value = 64998.75
$ws.Range("L44").Value = 64998.75
$ws.Range("N44").Value = -65974.75
$ws.Range("H45").Value = 2192.8125
$ws.Range("I45").Value = 1166.1666
$ws.Range("J45").Value = 2808.8
$ws.Range("K45").Value = 1166.1666
$ws.Range("L45").Value = 2808.8
$ws.Range("M45").Value = -789.1666
$ws.Range("N45").Value = -3562.8
$ws.Range("H55").Value = 49997.5
$ws.Range("J55").Value = 49997.5
$ws.Range("L55").Value = 49997.5
$ws.Range("N55").Value = -50627.5
$ws.Range("H102").Value = 2164.4
$ws.Range("I102").Value = 2170.2222
$ws.Range("J102").Value = 2112
$ws.Range("K102").Value = 2170.2222
$ws.Range("L102").Value = 2112
$ws.Range("M102").Value = -548.2222000000002
$ws.Range("N102").Value = -5356

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 46332
$ws.Range("J60").Value = 46332
$ws.Range("L60").Value = 46332
$ws.Range("N60").Value = -47530
$ws.Range("H99").Value = 5189
$ws.Range("I99").Value = 5111.25
$ws.Range("K99").Value = 5111.25
$ws.Range("M99").Value = -3613.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5425.8125
$ws.Range("I62").Value = 4661.4
$ws.Range("K62").Value = 4661.4
$ws.Range("M62").Value = -4037.4
$ws.Range("H65").Value = 5425.8125
$ws.Range("I65").Value = 4661.4
$ws.Range("K65").Value = 23307
$ws.Range("M65").Value = -20187
$ws.Range("H112").Value = 100702
$ws.Range("J112").Value = 100702
$ws.Range("L112").Value = 100702
$ws.Range("N112").Value = -103656

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1558.2858
$ws.Range("I2").Value = 2624.5
$ws.Range("K2").Value = 15747
$ws.Range("M2").Value = -15634
$ws.Range("H4").Value = 1046325.3
$ws.Range("I4").Value = 2092693.4
$ws.Range("K4").Value = 6278080.199999999
$ws.Range("M4").Value = -6277968.199999999
$ws.Range("H107").Value = 783.4706
$ws.Range("I107").Value = 434
$ws.Range("J107").Value = 1094.1111
$ws.Range("K107").Value = 1302
$ws.Range("L107").Value = 3282.3333
$ws.Range("M107").Value = 618
$ws.Range("N107").Value = -7122.3333
$ws.Range("H112").Value = 12812.714
$ws.Range("I112").Value = 4329.5
$ws.Range("J112").Value = 16206
$ws.Range("K112").Value = 12988.5
$ws.Range("L112").Value = 48618
$ws.Range("M112").Value = -11880.5
$ws.Range("N112").Value = -50834
$ws.Range("H113").Value = 1092.8334
$ws.Range("J113").Value = 1116.1428
$ws.Range("L113").Value = 3348.4284
$ws.Range("N113").Value = -7688.428400000001
$ws.Range("H115").Value = 733
$ws.Range("J115").Value = 999.5
$ws.Range("L115").Value = 2998.5
$ws.Range("N115").Value = -5348.5
$ws.Range("H131").Value = 4681866
$ws.Range("I131").Value = 11112118
$ws.Range("K131").Value = 33336354
$ws.Range("M131").Value = -33331314

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6224.3
$ws.Range("I70").Value = 5997
$ws.Range("J70").Value = 6754.6665
$ws.Range("K70").Value = 5997
$ws.Range("L70").Value = 6754.6665
$ws.Range("M70").Value = -5727
$ws.Range("N70").Value = -7294.6665
$ws.Range("H73").Value = 6224.3
$ws.Range("I73").Value = 5997
$ws.Range("J73").Value = 6754.6665
$ws.Range("K73").Value = 5997
$ws.Range("L73").Value = 6754.6665
$ws.Range("M73").Value = -5061
$ws.Range("N73").Value = -8626.666499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2032.5714
$ws.Range("I46").Value = 1023.4286
$ws.Range("K46").Value = 1023.4286
$ws.Range("M46").Value = -835.4286
$ws.Range("H55").Value = 618.5
$ws.Range("I55").Value = 438.2
$ws.Range("J55").Value = 1069.25
$ws.Range("K55").Value = 438.2
$ws.Range("L55").Value = 1069.25
$ws.Range("M55").Value = -265.2
$ws.Range("N55").Value = -1415.25
$ws.Range("H68").Value = 2204.389
$ws.Range("I68").Value = 2271.5
$ws.Range("J68").Value = 1969.5
$ws.Range("K68").Value = 2271.5
$ws.Range("L68").Value = 1969.5
$ws.Range("M68").Value = -1522.5
$ws.Range("N68").Value = -3467.5
$ws.Range("H71").Value = 2204.389
$ws.Range("I71").Value = 2271.5
$ws.Range("J71").Value = 1969.5
$ws.Range("K71").Value = 11357.5
$ws.Range("L71").Value = 9847.5
$ws.Range("M71").Value = -7613.5
$ws.Range("N71").Value = -17335.5
$ws.Range("H82").Value = 3366.6191
$ws.Range("J82").Value = 5129.9
$ws.Range("L82").Value = 5129.9
$ws.Range("N82").Value = -5851.9
$ws.Range("H85").Value = 3366.6191
$ws.Range("J85").Value = 5129.9
$ws.Range("L85").Value = 5129.9
$ws.Range("N85").Value = -7625.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3000
$ws.Range("I4").Value = 3000
$ws.Range("K4").Value = 3000
$ws.Range("M4").Value = -2887
$ws.Range("H81").Value = 5531.6665
$ws.Range("I81").Value = 6661.364
$ws.Range("K81").Value = 13322.728
$ws.Range("M81").Value = -12261.728
$ws.Range("H84").Value = 5531.6665
$ws.Range("I84").Value = 6661.364
$ws.Range("K84").Value = 66613.64
$ws.Range("M84").Value = -61309.64
$ws.Range("H122").Value = 2808.1667
$ws.Range("J122").Value = 2221
$ws.Range("L122").Value = 6663
$ws.Range("N122").Value = -11563
